$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the Model value in row 4 (Ground row) from "Box" to "Sphere"
$ws.Range("H4").Value = "Sphere"

# Update the selected cell/range shown in the sheet view
$ws.Range("D8").Select()
